# Update computed price/profit columns (H-N) across all 8 leve-profit sheets
# per scheduled data refresh.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 139.88889
$ws.Range("I9").Value = 82.5
$ws.Range("K9").Value = 82.5
$ws.Range("M9").Value = 86.5
$ws.Range("H10").Value = 268
$ws.Range("I10").Value = 268
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 268
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 25
$ws.Range("N10").ClearContents()
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H113").Value = 6575.5884
$ws.Range("I113").Value = 4050.8333
$ws.Range("J113").Value = 7952.727
$ws.Range("K113").Value = 4050.8333
$ws.Range("L113").Value = 7952.727
$ws.Range("M113").Value = -796.8332999999998
$ws.Range("N113").Value = -14460.727
$ws.Range("H135").Value = 3460.0625
$ws.Range("I135").Value = 3308.9285
$ws.Range("K135").Value = 29780.3565
$ws.Range("M135").Value = -27245.3565
$ws.Range("H138").Value = 3630.7102
$ws.Range("I138").Value = 2430
$ws.Range("J138").Value = 3810.8167
$ws.Range("K138").Value = 7290
$ws.Range("L138").Value = 11432.4501
$ws.Range("M138").Value = -2150
$ws.Range("N138").Value = -21712.4501

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 1099.6666
$ws.Range("J19").Value = 1799
$ws.Range("L19").Value = 1799
$ws.Range("N19").Value = -2257
$ws.Range("H32").Value = 8308.4375
$ws.Range("I32").Value = 8308.4375
$ws.Range("K32").Value = 8308.4375
$ws.Range("M32").Value = -8021.4375
$ws.Range("H45").Value = 2972.2727
$ws.Range("I45").Value = 2186.7778
$ws.Range("K45").Value = 2186.7778
$ws.Range("M45").Value = -1809.7778
$ws.Range("H61").Value = 5935.5293
$ws.Range("I61").Value = 4493.125
$ws.Range("K61").Value = 4493.125
$ws.Range("M61").Value = -4281.125
$ws.Range("H132").Value = 2739.457
$ws.Range("I132").Value = 2217.6072
$ws.Range("J132").Value = 4826.857
$ws.Range("K132").Value = 6652.821599999999
$ws.Range("L132").Value = 14480.571
$ws.Range("M132").Value = -4122.821599999999
$ws.Range("N132").Value = -19540.571
$ws.Range("H136").Value = 5935.5293
$ws.Range("I136").Value = 4493.125
$ws.Range("K136").Value = 13479.375
$ws.Range("M136").Value = -10929.375

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 23642.572
$ws.Range("I82").Value = 10916.5
$ws.Range("J82").Value = 99999
$ws.Range("K82").Value = 10916.5
$ws.Range("L82").Value = 99999
$ws.Range("M82").Value = -10533.5
$ws.Range("N82").Value = -100765
$ws.Range("H85").Value = 23642.572
$ws.Range("I85").Value = 10916.5
$ws.Range("J85").Value = 99999
$ws.Range("K85").Value = 10916.5
$ws.Range("L85").Value = 99999
$ws.Range("M85").Value = -9590.5
$ws.Range("N85").Value = -102651
$ws.Range("H86").Value = 3155.5386
$ws.Range("I86").Value = 2451.8096
$ws.Range("J86").Value = 6111.2
$ws.Range("K86").Value = 2451.8096
$ws.Range("L86").Value = 6111.2
$ws.Range("M86").Value = -1328.8096
$ws.Range("N86").Value = -8357.200000000001
$ws.Range("H89").Value = 3155.5386
$ws.Range("I89").Value = 2451.8096
$ws.Range("J89").Value = 6111.2
$ws.Range("K89").Value = 12259.048
$ws.Range("L89").Value = 30556
$ws.Range("M89").Value = -6643.048000000001
$ws.Range("N89").Value = -41788
$ws.Range("H107").Value = 1679.421
$ws.Range("I107").Value = 1761.6875
$ws.Range("K107").Value = 1761.6875
$ws.Range("M107").Value = 158.3125
$ws.Range("H134").Value = 4401.077
$ws.Range("I134").Value = 3114.4285
$ws.Range("J134").Value = 5902.1665
$ws.Range("K134").Value = 9343.2855
$ws.Range("L134").Value = 17706.4995
$ws.Range("M134").Value = -6808.2855
$ws.Range("N134").Value = -22776.4995

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 4066.3333
$ws.Range("I6").Value = 4066.3333
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 4066.3333
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -3953.3333
$ws.Range("N6").ClearContents()
$ws.Range("H10").Value = 7499.5
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 7499.5
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 7499.5
$ws.Range("N10").Value = -7777.5
$ws.Range("M10").ClearContents()
$ws.Range("H13").Value = 4333.3335
$ws.Range("J13").Value = 4333.3335
$ws.Range("L13").Value = 4333.3335
$ws.Range("N13").Value = -4611.3335
$ws.Range("H22").Value = 6610.636
$ws.Range("I22").Value = 5539.6
$ws.Range("J22").Value = 7503.1665
$ws.Range("K22").Value = 5539.6
$ws.Range("L22").Value = 7503.1665
$ws.Range("M22").Value = -5189.6
$ws.Range("N22").Value = -8203.166499999999
$ws.Range("H31").Value = 45717.44
$ws.Range("I31").Value = 2258.75
$ws.Range("K31").Value = 2258.75
$ws.Range("M31").Value = -1963.75
$ws.Range("H34").Value = 45717.44
$ws.Range("I34").Value = 2258.75
$ws.Range("K34").Value = 2258.75
$ws.Range("M34").Value = -2056.75
$ws.Range("H99").Value = 2883.5217
$ws.Range("I99").Value = 3049.3333
$ws.Range("K99").Value = 3049.3333
$ws.Range("M99").Value = -1551.3333
$ws.Range("H107").Value = 1404.2858
$ws.Range("I107").Value = 639.619
$ws.Range("J107").Value = 3698.2856
$ws.Range("K107").Value = 639.619
$ws.Range("L107").Value = 3698.2856
$ws.Range("M107").Value = 1280.381
$ws.Range("N107").Value = -7538.2856
$ws.Range("H126").Value = 2883.5217
$ws.Range("I126").Value = 3049.3333
$ws.Range("K126").Value = 9147.999899999999
$ws.Range("M126").Value = -6677.999899999999
$ws.Range("H135").Value = 59116
$ws.Range("I135").Value = 58448
$ws.Range("J135").Value = 59199.5
$ws.Range("K135").Value = 58448
$ws.Range("L135").Value = 59199.5
$ws.Range("M135").Value = -53378
$ws.Range("N135").Value = -69339.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 92884.17999999999
$ws.Range("I11").Value = 168455
$ws.Range("K11").Value = 505365
$ws.Range("M11").Value = -505225
$ws.Range("H34").Value = 8522735
$ws.Range("I34").Value = 25929326
$ws.Range("J34").Value = 4171087.5
$ws.Range("K34").Value = 77787978
$ws.Range("L34").Value = 12513262.5
$ws.Range("M34").Value = -77787894
$ws.Range("N34").Value = -12513430.5
$ws.Range("H55").Value = 2215.5
$ws.Range("J55").Value = 5500
$ws.Range("L55").Value = 16500
$ws.Range("N55").Value = -16854
$ws.Range("H140").Value = 4743.826
$ws.Range("I140").Value = 3132.2307
$ws.Range("K140").Value = 9396.6921
$ws.Range("M140").Value = -4216.6921

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5388.0835
$ws.Range("I132").Value = 3592.647
$ws.Range("J132").Value = 9748.429
$ws.Range("K132").Value = 10777.941
$ws.Range("L132").Value = 29245.287
$ws.Range("M132").Value = -8247.940999999999
$ws.Range("N132").Value = -34305.287

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6697.92
$ws.Range("I40").Value = 5181.2104
$ws.Range("K40").Value = 5181.2104
$ws.Range("M40").Value = -5045.2104
$ws.Range("H46").Value = 2890
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5376
$ws.Range("H100").Value = 4008.4666
$ws.Range("I100").Value = 2801.5833
$ws.Range("K100").Value = 2801.5833
$ws.Range("M100").Value = -2260.5833
$ws.Range("H132").Value = 5341.8667
$ws.Range("I132").Value = 3697.7827
$ws.Range("K132").Value = 11093.3481
$ws.Range("M132").Value = -8563.348100000001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1841.8572
$ws.Range("I132").Value = 1454.6757
$ws.Range("J132").Value = 4707
$ws.Range("K132").Value = 4364.0271
$ws.Range("L132").Value = 14121
$ws.Range("M132").Value = -1834.0271
$ws.Range("N132").Value = -19181
$ws.Range("H136").Value = 2239.9534
$ws.Range("I136").Value = 1239.7059
$ws.Range("K136").Value = 3719.1177
$ws.Range("M136").Value = -1169.1177
